$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.114.86"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").Value = "1.928.13"
$ws.Range("E3").Value = "  +2.18%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.17"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("E7").Value = "  +1.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3830"
$ws.Range("E8").Value = "  +0.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07767"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9819"
$ws.Range("E10").Value = "  +2.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.73"
$ws.Range("E11").Value = "  +3.44%  "
$ws.Range("D12").Value = "1.934.32"
$ws.Range("E12").Value = "  +2.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.710"
$ws.Range("E13").Value = "  +1.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.978"
$ws.Range("E14").Value = "  +0.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07051"
$ws.Range("E15").Value = "  +0.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.006"
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "84.44"
$ws.Range("E17").Value = "  +1.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009553"
$ws.Range("E18").Value = "  +1.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.77"
$ws.Range("E19").Value = "  +1.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("E20").Value = "  +0.31%  "
$ws.Range("D21").Value = "29.114.29"
$ws.Range("E21").Value = "  +1.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.356"
$ws.Range("E22").Value = "  +0.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.98"
$ws.Range("E23").Value = "  +1.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.073"
$ws.Range("E24").Value = "  -0.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.65"
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "19.11"
$ws.Range("E26").Value = "  +0.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.680"
$ws.Range("E27").Value = "  +1.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "118.15"
$ws.Range("E28").Value = "  +1.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.846"
$ws.Range("E29").Value = "  +2.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09354"
$ws.Range("E30").Value = "  +1.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.8619"
$ws.Range("E31").Value = "  +2.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.126"
$ws.Range("E32").Value = "  +1.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.249"
$ws.Range("E33").Value = "  +1.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.015"
$ws.Range("E34").Value = "  +0.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.163"
$ws.Range("E35").Value = "  +1.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05702"
$ws.Range("E36").Value = "  +0.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.215"
$ws.Range("E37").Value = "  +19.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.005"
$ws.Range("E38").Value = "  +0.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02053"
$ws.Range("E39").Value = "  +1.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.528"
$ws.Range("E40").Value = "  +1.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5530"
$ws.Range("E41").Value = "  +0.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1759"
$ws.Range("E42").Value = "  +0.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.376"
$ws.Range("E43").Value = "  +2.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.197"
$ws.Range("E44").Value = "  +6.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.000002747"
$ws.Range("E45").Value = "  -6.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5214"
$ws.Range("E46").Value = "  +1.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.33"
$ws.Range("E47").Value = "  +1.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06930"
$ws.Range("E48").Value = "  +1.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "110.58"
$ws.Range("E49").Value = "  -0.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.775"
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.004"
$ws.Range("E51").Value = "  +0.29%  "
